# Update crypto price symbols (column D) on the active sheet.
# Cells are stored as text (inlineStr) values, so we force the
# NumberFormat to Text before assigning the new values to keep the
# exact textual representation (e.g. trailing zeros) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "271.23"
    "D3"  = "22.80"
    "D4"  = "6.375"
    "D5"  = "0.06278"
    "D7"  = "6.711"
    "D8"  = "1.377"
    "D9"  = "0.8360"
    "D10" = "0.01377"
    "D11" = "0.1630"
    "D12" = "0.08416"
    "D13" = "0.03488"
    "D14" = "0.03143"
    "D15" = "0.09320"
    "D16" = "3.886"
    "D17" = "0.001715"
    "D18" = "0.04822"
    "D19" = "0.006203"
    "D21" = "0.003614"
    "D22" = "0.0001497"
    "D23" = "3.736"
    "D25" = "0.3404"
    "D26" = "0.1262"
    "D40" = "0.04687"
    "D41" = "0.006924"
    "D42" = "0.1173"
    "D44" = "0.01253"
    "D45" = "0.00006260"
    "D47" = "0.7967"
    "D48" = "0.09077"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
